# Generate Report for Archive
# Status moved on from "Ready for handoff" -> "In Translation" for the
# two tracked source files, across the Overview rollup and each locale
# report sheet. Updating the text also narrows the status columns (the
# new text is shorter), matching the regenerated report's column sizing.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$overviewRows = $overview.UsedRange.Rows.Count
for ($r = 2; $r -le $overviewRows; $r++) {
    if ($overview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value2 = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value2 = $newStatus
    }
}

# --- Locale sheets: single "Status" column (C) ---
foreach ($sheet in @($zhcn, $dede)) {
    $rows = $sheet.UsedRange.Rows.Count
    for ($r = 2; $r -le $rows; $r++) {
        if ($sheet.Cells.Item($r, 3).Value2 -eq $oldStatus) {
            $sheet.Cells.Item($r, 3).Value2 = $newStatus
        }
    }
}

# --- Re-size the status columns to fit the new (shorter) text ---
# Target "character" width ~13.41 (engine quantizes ColumnWidth to 1/6
# increments before storing, so 12.5 is the input that lands closest to
# that target once written back out as the sheet's <col width>).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
